# InputParam1.xlsx update:
#   - Row 43/44 swap which fitting parameter they represent
#     (g_PMCA now comes before nu_leakSR), keeping each parameter's own value.
#   - A new parameter row (g_leakNa) is appended as row 45.
#   - Selection/view is moved down near the newly-added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 43 now holds g_PMCA (previously held nu_leakSR).
$ws.Range("A43").Value = "g_PMCA"
$ws.Range("B43").Value = 5.37

# Row 44 now holds nu_leakSR (previously held g_PMCA).
$ws.Range("A44").Value = "nu_leakSR"
$ws.Range("B44").Value = 0.2

# New row 45: g_leakNa parameter.
$ws.Range("A45").Value = "g_leakNa"
$ws.Range("B45").Value = 0.002

# Update selection to reflect where editing continued (B46, just past the new row).
$ws.Range("B46").Select()
